# Applies the "Cleanedup weapons/armour" edit to the censers sheet:
#   - rewrites the L ("cost") column values for rows 2..60 with the
#     rebalanced cost curve
#   - gives the L column an integer ("0") number format, distinct from the
#     General format used before
#   - moves the active selection/top-left cell back to the top of the
#     sheet (A1) with L2:L60 selected, instead of AE2:AE60 scrolled to AZ1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "cost" values for L2:L60 (rows 2-60)
$costValues = @(
    9,13,19,28,40,59,87,127,186,273,
    399,585,856,1254,1837,2691,3941,5772,8454,12382,
    18135,26561,38902,56977,83451,122226,179017,262195,384022,562454,
    823793,1206561,1767179,2588282,3790904,5552314,8132146,11910675,17444863,25550461,
    37422253,54810166,80277216,117577302,172208537,252223684,369417149,541063502,792463789,1160674956,
    1699972128,2000000000,3646734750,5341157231,7822877869,11457707664,16781428411,24578768089,35999071473
)

for ($i = 0; $i -lt $costValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 12).Value = $costValues[$i]
}

# Give the cost column (L2:L60) a plain integer display format.
$ws.Range("L2:L60").NumberFormat = "0"

# Restore the sheet view: scroll back to the top-left (A1) and select L2:L60
# (previously the view was scrolled to AZ1 with AE2:AE60 selected).
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("L2:L60").Select()
$excel.ActiveWindow.DisplayGridlines = $true
$excel.ActiveWindow.DisplayHeadings = $true
$excel.ActiveWindow.DisplayZeros = $true
